# Apply the two content edits described in the commit:
# 1. "Ran home-based ... gaming PC" + bookmark + "s." -> merge into a single
#    run reading "...gaming PCs." with the old _GoBack bookmark removed.
# 2. Insert " MTA SQL certified" after "...MTA Network certified," and move
#    the _GoBack bookmark to sit right after the newly typed text (i.e.
#    right before " and familiarity with video editing tools such as Movie
#    Maker"), matching where Word leaves the caret after typing.

$d = $word.ActiveDocument

# --- Edit 1: "gaming PC" + bookmark + "s." -> "gaming PCs." ---------------------
# A Find/Replace that spans across the existing _GoBack bookmark collapses the
# bookmark and merges the surrounding runs, exactly like Word does when text
# spanning a bookmark is edited/retyped.
$d.Content.Find.Execute(
    "Ran home-based custom computer and repair service in Benzie County specializing in high-performance gaming PCs.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ran home-based custom computer and repair service in Benzie County specializing in high-performance gaming PCs.",
    2) | Out-Null

# --- Edit 2: add " MTA SQL certified" and relocate the _GoBack bookmark --------
$found = $d.Content
$found.Find.Execute("MTA Network certified,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$insertionPoint = $found.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.InsertAfter(" MTA SQL certified")

# After InsertAfter, $insertionPoint's End marks the spot right after the
# freshly typed text -- that's where Word leaves (and re-stamps) _GoBack.
$caret = $d.Range($insertionPoint.End, $insertionPoint.End)
$d.Bookmarks.Add("_GoBack", $caret) | Out-Null
